# Insert a new price-report row at row 6 (weekly update), pushing the
# existing rows 6-24 down to 7-25 and appending one more row of history.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 6..24 down one slot, leaving a blank (but formatted) row 6.
$ws.Rows("6:6").Insert()

# Populate the new row 6 with this week's report.
$ws.Cells.Item(6, 1).Value  = 3
$ws.Cells.Item(6, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(6, 3).Value  = "Coquimbo"
$ws.Cells.Item(6, 4).Value  = 44469
$ws.Cells.Item(6, 5).Value  = 5
$ws.Cells.Item(6, 6).Value  = 100112022
$ws.Cells.Item(6, 7).Value  = "Arveja Verde"
$ws.Cells.Item(6, 8).Value  = "Perfection"
$ws.Cells.Item(6, 9).Value  = "Primera"
$ws.Cells.Item(6, 10).Value = 73
$ws.Cells.Item(6, 11).Value = 28000
$ws.Cells.Item(6, 12).Value = 29000
$ws.Cells.Item(6, 13).Value = 28521
$ws.Cells.Item(6, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(6, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(6, 16).Value = 1141
$ws.Cells.Item(6, 17).Value = 25
$ws.Cells.Item(6, 18).Value = "Hortaliza"
